$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1 and H1
$ws.Range("G1").Value = "Latest Update User"
$ws.Range("H1").Value = "Latest Update Date"

# Copy the style/format from F1 (existing header cell) to the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active cell / selection to H2, scrolled so column E is the leftmost visible column
$ws.Range("H2").Select()
$excel.ActiveWindow.ScrollColumn = 5
